# Apply cryptos list update (cell value changes) to match target diff
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Cells whose new value is a clean decimal-looking number but must remain text ---
# (matches original column D formatting, which stores every price as text)
$textForceCells = @('D4', 'D5', 'D6', 'D8', 'D9', 'D13', 'D15', 'D16', 'D19', 'D24', 'D25', 'D26', 'D28', 'D29', 'D32', 'D33', 'D34', 'D38', 'D39', 'D40', 'D41', 'D42', 'D44', 'D48', 'D49', 'D50')
foreach ($cellRef in $textForceCells) {
    $ws.Range($cellRef).NumberFormat = "@"
}

# --- Set updated values ---
$ws.Range('D4').Value = '0.9993'
$ws.Range('D5').Value = '241.11'
$ws.Range('D6').Value = '0.6713'
$ws.Range('D8').Value = '0.07436'
$ws.Range('D9').Value = '0.2939'
$ws.Range('D13').Value = '5.008'
$ws.Range('D15').Value = '85.98'
$ws.Range('D16').Value = '6.159'
$ws.Range('D19').Value = '228.36'
$ws.Range('D24').Value = '160.85'
$ws.Range('D25').Value = '8.704'
$ws.Range('D26').Value = '0.1403'
$ws.Range('D28').Value = '1.512'
$ws.Range('D29').Value = '4.159'
$ws.Range('D32').Value = '0.05292'
$ws.Range('D33').Value = '1.876'
$ws.Range('D34').Value = '0.7521'
$ws.Range('D38').Value = '0.01806'
$ws.Range('D39').Value = '2.729'
$ws.Range('D40').Value = '0.9215'
$ws.Range('D41').Value = '5.966'
$ws.Range('D42').Value = '0.08364'
$ws.Range('D44').Value = '102.18'
$ws.Range('D48').Value = '0.00000000121'
$ws.Range('D49').Value = '63.75'
$ws.Range('D50').Value = '9.150'
$ws.Range('D2').Value = '29.281.92'
$ws.Range('E2').Value = '  +0.48%  '
$ws.Range('E3').Value = '  +0.10%  '
$ws.Range('E5').Value = '  -0.13%  '
$ws.Range('E6').Value = '  -2.26%  '
$ws.Range('E7').Value = '  +0.07%  '
$ws.Range('E8').Value = '  -0.20%  '
$ws.Range('E9').Value = '  -2.57%  '
$ws.Range('E10').Value = '  -1.00%  '
$ws.Range('E11').Value = '  +0.93%  '
$ws.Range('D12').Value = '1.833.85'
$ws.Range('E12').Value = '  -0.30%  '
$ws.Range('E13').Value = '  -1.03%  '
$ws.Range('E14').Value = '  -1.67%  '
$ws.Range('E15').Value = '  -1.82%  '
$ws.Range('E16').Value = '  -0.21%  '
$ws.Range('D17').Value = '29.241.03'
$ws.Range('E17').Value = '  +0.38%  '
$ws.Range('E18').Value = '  +2.04%  '
$ws.Range('E19').Value = '  -0.01%  '
$ws.Range('E20').Value = '  -0.22%  '
$ws.Range('E21').Value = '  +0.10%  '
$ws.Range('E22').Value = '  -3.10%  '
$ws.Range('E24').Value = '  +0.64%  '
$ws.Range('E25').Value = '  -0.65%  '
$ws.Range('E26').Value = '  -3.37%  '
$ws.Range('E27').Value = '  -0.34%  '
$ws.Range('E28').Value = '  +0.08%  '
$ws.Range('E29').Value = '  -2.68%  '
$ws.Range('E30').Value = '  -1.76%  '
$ws.Range('E31').Value = '  +0.03%  '
$ws.Range('E32').Value = '  +0.85%  '
$ws.Range('E33').Value = '  +1.36%  '
$ws.Range('E34').Value = '  -0.89%  '
$ws.Range('E35').Value = '  +0.28%  '
$ws.Range('E36').Value = '  -0.23%  '
$ws.Range('D37').Value = '1.320.84'
$ws.Range('E37').Value = '  +1.37%  '
$ws.Range('E38').Value = '  -1.55%  '
$ws.Range('E39').Value = '  +0.17%  '
$ws.Range('E40').Value = '  -0.86%  '
$ws.Range('E41').Value = '  +0.86%  '
$ws.Range('E42').Value = '  +12.79%  '
$ws.Range('E43').Value = '  +0.28%  '
$ws.Range('E44').Value = '  -2.58%  '
$ws.Range('D45').Value = '1.975.82'
$ws.Range('E45').Value = '  -0.27%  '
$ws.Range('E46').Value = '  -0.66%  '
$ws.Range('E47').Value = '  +0.29%  '
$ws.Range('B48').Value = 'BabyDogeCoin'
$ws.Range('C48').Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range('E48').Value = '  -1.77%  '
$ws.Range('B49').Value = 'Aave'
$ws.Range('C49').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('E49').Value = '  -1.97%  '
$ws.Range('E50').Value = '  -4.00%  '
$ws.Range('E51').Value = '  +0.03%  '
